# Add a "Greece" market tab, cloned from the "Croatia" tab, and fill in
# the Greece-specific values (mirrors commit "Test data for Greece Market").

$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Duplicate the Croatia sheet immediately after itself, then rename it.
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Update the market-specific values on the new sheet.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3189"

# Leave the Croatia source tab with its whole sheet highlighted (as it is
# left after copying its contents out), then make Greece the active tab
# with B4 selected.
$croatia.Activate()
$croatia.Cells.Select()

$greece.Activate()
$greece.Range("B4").Select()
